$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("BC4").Value = 151
$ws.Range("O6").Value = 1.14
$ws.Range("P6").Value = 5.5
$ws.Range("G7").Value = 1.8
$ws.Range("H7").Value = 3.4
$ws.Range("I7").Value = 4.75
$ws.Range("L7").Value = 5
$ws.Range("AG7").Value = 12
$ws.Range("AH7").Value = 23
$ws.Range("AM7").Value = 351
$ws.Range("AP7").Value = 21
$ws.Range("AX7").Value = 26
$ws.Range("BA7").Value = 126
$ws.Range("G8").Value = 2.2
$ws.Range("H8").Value = 3.3
$ws.Range("I8").Value = 3.2
$ws.Range("J8").Value = 3
$ws.Range("K8").Value = 2
$ws.Range("O8").Value = 1.4
$ws.Range("P8").Value = 2.75
$ws.Range("Q8").Value = 2.25
$ws.Range("R8").Value = 1.62
$ws.Range("S8").Value = 1.5
$ws.Range("T8").Value = 2.5
$ws.Range("U8").Value = 2
$ws.Range("V8").Value = 1.73
$ws.Range("W8").Value = 6.5
$ws.Range("Z8").Value = 21
$ws.Range("AA8").Value = 21
$ws.Range("AC8").Value = 8
$ws.Range("AF8").Value = 67
$ws.Range("AG8").Value = 8.5
$ws.Range("AO8").Value = 13
$ws.Range("AP8").Value = 26
$ws.Range("AT8").Value = 2.5
$ws.Range("AV8").Value = 67
$ws.Range("O9").Value = 1.29
$ws.Range("P9").Value = 3.75
$ws.Range("Q9").Value = 1.97
$ws.Range("R9").Value = 1.93
$ws.Range("Q10").Value = 2.35
$ws.Range("R10").Value = 1.57
$ws.Range("N13").Value = 8
$ws.Range("G16").Value = 1.75
$ws.Range("J16").Value = 2.5
$ws.Range("M16").Value = 1.11
$ws.Range("N16").Value = 6.5
$ws.Range("AC16").Value = 6.5
$ws.Range("AQ16").Value = 34
$ws.Range("H17").Value = 3.2
$ws.Range("I17").Value = 3.9
$ws.Range("Z17").Value = 17
$ws.Range("AA17").Value = 19
$ws.Range("AD17").Value = 6
$ws.Range("AG17").Value = 9
$ws.Range("AH17").Value = 19
$ws.Range("AI17").Value = 15
$ws.Range("AK17").Value = 41
$ws.Range("G20").Value = 30
$ws.Range("I20").Value = 1.06
$ws.Range("J20").Value = 20
$ws.Range("K20").Value = 3.6
$ws.Range("L20").Value = 1.28
$ws.Range("P20").Value = 6.8
$ws.Range("R20").Value = 3.6
$ws.Range("W20").Value = 150
$ws.Range("X20").Value = 800
$ws.Range("Y20").Value = 150
$ws.Range("AB20").Value = 400
$ws.Range("AC20").Value = 25
$ws.Range("AD20").Value = 24
$ws.Range("AE20").Value = 50
$ws.Range("AF20").Value = 200
$ws.Range("AH20").Value = 7.5
$ws.Range("AI20").Value = 13.5
$ws.Range("AK20").Value = 11.5
$ws.Range("AL20").Value = 40
$ws.Range("AN20").Value = 30
$ws.Range("AO20").Value = 250
$ws.Range("AP20").Value = 110
$ws.Range("AU20").Value = 11.75
$ws.Range("AV20").Value = 90
$ws.Range("AZ20").Value = 6.8
